$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.758.17"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "2.216.51"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "292.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "86.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.01%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0785"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("E13").Value = "  +1.63%  "
$ws.Range("E14").Value = "  +1.96%  "
$ws.Range("D15").Value = "2.561.46"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").Value = "2.222.70"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("E18").Value = "  +3.47%  "
$ws.Range("D19").Value = "39.731.86"
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("D20").Value = "0.0₃0881"
$ws.Range("E20").Value = "  +1.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.12%  "
$ws.Range("E22").Value = "  +2.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.58%  "
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("E26").Value = "  +2.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("E29").Value = "  +1.70%  "
$ws.Range("E30").Value = "  +2.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.49%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.07%  "
$ws.Range("E35").Value = "  +4.29%  "
$ws.Range("E36").Value = "  +2.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.27%  "
$ws.Range("E38").Value = "  +1.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0989"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.94%  "
$ws.Range("E41").Value = "  +4.75%  "
$ws.Range("E42").Value = "  +5.83%  "
$ws.Range("D43").Value = "2.068.82"
$ws.Range("E43").Value = "  +9.56%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.67%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0267"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.52%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.52%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("D49").Value = "2.433.10"
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "88.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.80%  "
